$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1, "sum") onto the
# new header cell H1 so the new column reuses the same bold/bordered header
# style instead of minting a near-duplicate style entry.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New header label for the "Save" column.
$ws.Range("H1").Value = "Save"

# New "Save" column values (plain, unstyled numeric cells), same as the
# other data columns.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
